# Applies the committed changes:
#  - Row 4: round Q4/R4 (Ost/Nord coordinates) to nearest integer
#  - Rows 5 and 6: the two fungi observations were swapped (A, B, E, F, G, H,
#    Z, AB) and their Q/R coordinates were rounded to nearest integer.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: round Ost/Nord coordinates ---
$ws.Range("Q4").Value = 690281
$ws.Range("R4").Value = 7126404

# --- Capture current (pre-swap) values for rows 5 and 6 ---
# NOTE: use .Value2 for reads (.Value reads have been observed to yield the
# raw property-accessor object rather than the cell's content in this host).
$A5 = $ws.Range("A5").Value2
$B5 = $ws.Range("B5").Value2
$E5 = $ws.Range("E5").Value2
$F5 = $ws.Range("F5").Value2
$G5 = $ws.Range("G5").Value2
$H5 = $ws.Range("H5").Value2
$Z5 = $ws.Range("Z5").Value2
$AB5 = $ws.Range("AB5").Value2

$A6 = $ws.Range("A6").Value2
$B6 = $ws.Range("B6").Value2
$E6 = $ws.Range("E6").Value2
$F6 = $ws.Range("F6").Value2
$G6 = $ws.Range("G6").Value2
$H6 = $ws.Range("H6").Value2
$Z6 = $ws.Range("Z6").Value2
$AB6 = $ws.Range("AB6").Value2

# --- Write row 5 with row 6's former values ---
$ws.Range("A5").Value = $A6
$ws.Range("B5").Value = $B6
$ws.Range("E5").Value = $E6
$ws.Range("F5").Value = $F6
$ws.Range("G5").Value = $G6
$ws.Range("H5").Value = $H6
$ws.Range("Z5").Value = $Z6
$ws.Range("AB5").Value = $AB6

# --- Write row 6 with row 5's former values ---
$ws.Range("A6").Value = $A5
$ws.Range("B6").Value = $B5
$ws.Range("E6").Value = $E5
$ws.Range("F6").Value = $F5
$ws.Range("G6").Value = $G5
$ws.Range("H6").Value = $H5
$ws.Range("Z6").Value = $Z5
$ws.Range("AB6").Value = $AB5

# --- Q/R coordinates for rows 5 and 6: swap and round ---
$ws.Range("Q5").Value = 690408
$ws.Range("R5").Value = 7125570
$ws.Range("Q6").Value = 690447
$ws.Range("R6").Value = 7125629
